# Natmi following Dr Hou advice
# Update LR-pair results for Efna5-Ephb6: add an "ECs" sending-cluster group
# and refresh the computed specificity/weight values for all sending/target
# cluster combinations.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Efna5", "Ephb6", "FAPs", 2, 0.6666666666666666, 0.4223226666666666, 1.266968, 0.1676547342089268, 0.1676547342089269, 3, 1, 0.7341896666666666, 2.202569, 0.4912907638668469, 0.4912907638668469, 0.3100649378657777, 2.790584440792, 0.08236722243539686, 0.08236722243539688),
    @("ECs", "Efna5", "Ephb6", "sCs", 2, 0.6666666666666666, 0.4223226666666666, 1.266968, 0.1676547342089268, 0.1676547342089269, 3, 1, 0.7602199999999999, 2.28066, 0.5087092361331531, 0.5087092361331531, 0.3210581376533332, 2.889523238879999, 0.08528751177352999, 0.08528751177353),
    @("FAPs", "Efna5", "Ephb6", "FAPs", 3, 1, 1.874986333333333, 5.624959, 0.7443368783435028, 0.7443368783435029, 3, 1, 0.7341896666666666, 2.202569, 0.4912907638668469, 0.4912907638668469, 1.376595591074556, 12.389360319671, 0.3656858335356438, 0.3656858335356438),
    @("FAPs", "Efna5", "Ephb6", "sCs", 3, 1, 1.874986333333333, 5.624959, 0.7443368783435028, 0.7443368783435029, 3, 1, 0.7602199999999999, 2.28066, 0.5087092361331531, 0.5087092361331531, 1.425402110326667, 12.82861899294, 0.378651044807859, 0.3786510448078591),
    @("sCs", "Efna5", "Ephb6", "FAPs", 3, 1, 0.2216933333333333, 0.66508, 0.08800838744757017, 0.08800838744757018, 3, 1, 0.7341896666666666, 2.202569, 0.4912907638668469, 0.4912907638668469, 0.1627649545022222, 1.46488459052, 0.04323770789580617, 0.04323770789580617),
    @("sCs", "Efna5", "Ephb6", "sCs", 3, 1, 0.2216933333333333, 0.66508, 0.08800838744757017, 0.08800838744757018, 3, 1, 0.7602199999999999, 2.28066, 0.5087092361331531, 0.5087092361331531, 0.1685357058666666, 1.5168213528, 0.044770679551764, 0.04477067955176401),
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowIndex = $i + 2
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $rowValues.Count; $j++) {
        $ws.Cells.Item($rowIndex, $j + 1).Value = $rowValues[$j]
    }
}
